# Main.xlsx - DesignFirst project save
# 1) Mark the used columns (A:L, i.e. the whole "custom width" block plus the
#    trailing default-width run) as collapsed, matching the outline/grouping
#    state recorded for this sheet.
# 2) D10 changes from 21 to 100 (matching the value already present in C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = $ws.Columns("A:L")
$cols.Collapsed = $true

$ws.Range("D10").Value = 100.0
